$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column R (2021) mirroring the styles of column Q (2020).
# Copying the existing Q3:Q34 range seeds column R with identical cell
# styles/number formats (and bumps dimension + row "spans" automatically),
# after which the actual 2021 figures are written on top.
$ws.Range("Q3:Q34").Copy($ws.Range("R3:R34"))

$values = @{
    4  = 2021
    5  = 109
    6  = 74
    7  = 35
    8  = 36
    9  = 35
    10 = 1
    11 = 15
    12 = 8
    13 = 7
    14 = 12
    15 = 7
    16 = 5
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 17
    21 = 8
    22 = 9
    23 = 9
    24 = 7
    25 = 2
    26 = 20
    27 = 9
    28 = 11
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

foreach ($row in $values.Keys) {
    $ws.Range("R$row").Value2 = $values[$row]
}

# Match the saved selection state from the edit (activeCell moved to R1).
$ws.Range("R1").Select()
